# Fix driver path for UI Engine.
# Re-running the UI test cases (Home.01, Search.01, Search.02) failed because the
# chromedriver executable had the wrong permissions. Record the new run results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests1")

$runTime = "2017-10-22 16:51:52"
$failMsg = "[Line 0] Message: '' executable may have wrong permissions. Please see https://sites.google.com/a/chromium.org/chromedriver/home`n"

# AT.UI.Home.01 (row 2)
$ws.Range("J2").Value = "fail"
$ws.Range("K2").Value = $runTime
$ws.Range("L2").Value = $failMsg

# AT.UI.Search.01 (row 6)
$ws.Range("J6").Value = "fail"
$ws.Range("K6").Value = $runTime
$ws.Range("L6").Value = $failMsg

# AT.UI.Search.02 (row 11)
$ws.Range("J11").Value = "fail"
$ws.Range("K11").Value = $runTime
$ws.Range("L11").Value = $failMsg

# AT.HTTP.Get.03 (row 42) - also re-run (still passes), only timestamp updates
$ws.Range("K42").Value = $runTime
